# Close defects 64-67 and 80-87
$d = $word.ActiveDocument

# --- 1) Add the "Date Closed" value "17/5/2023" -----------------------
# Locate the "Date Closed: " label (including its trailing space run),
# collapse to the end of that text, and type the date there so it lands
# in a new run right after the existing space run.
$rng = $d.Content
$found = $rng.Find.Execute("Date Closed: ", $false, $false, $false, $false,
                            $false, $true, 1, $false, "", 0)

if ($found) {
    $rng.Collapse(0)
    $rng.Select()
    $sel = $word.Selection
    $sel.InsertAfter("17/5/2023")
    $sel.Font.Bold = $true
    $sel.Font.BoldBi = $true
    $sel.Font.Size = 20
    $sel.Font.SizeBi = 20
}

# --- 2) Mark the "Screenshot for the Defect:" run as NoProof ----------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Screenshot for the Defect:", $false, $false,
                              $false, $false, $false, $true, 1, $false,
                              "", 0)
if ($found2) {
    $rng2.NoProofing = 1
}
